$wb = $excel.ActiveWorkbook

# Update the status text "Ready for handoff" -> "In Translation" wherever it occurs
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2
        if ("Ready for handoff" -eq $v) {
            $cell.Value = "In Translation"
        }
    }
}

# Autofit the affected (Status) columns on each sheet so widths recalc based on new text
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$wsOverview.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
